# Script 1 - atualização automática de dados (Execução: 23)
# Updates the numeric data in columns B:D (rows 2-13) of the active sheet
# with freshly recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2012, -6.94354713543871,  7.990168510511264,   3.559831809881109),
    @(2013, -3.033464260293317, 4.333352843923399,  -5.227846237893674),
    @(2014, -0.2234743598115374,4.093761345842939,   1.620256385538821),
    @(2015,  1.260711594855279,-1.223239834604506,   8.070903554328556),
    @(2016, -5.160995246877953,-3.632675769107285,   0.1732976787794716),
    @(2017, -3.381558148626762, 0.5958395557011942,  0.7220832816449141),
    @(2018, -3.186190652963306,-1.091206871444617,  -1.955526309528577),
    @(2019,  2.964647058339054, 1.042552283805143,   9.570679981139186),
    @(2020,-13.99768529350024, -3.93685177458396,  -12.15959139070785),
    @(2021,-11.0410856605323,  15.37783447774446,  -14.1120775080652),
    @(2022, -4.78117406122619, 14.54741383364233,  -12.90242013598057),
    @(2023, -5.978844308965603, 7.560548501577813,  -7.493641516982841)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $row++
}
